# Refresh the cached regression-table figures on the "multiple_loans" sheet
# (stabilized cleaning pipeline + internal-validity re-run for the number of
# pawns balance). The table's data cells are formulas that read cached
# results from an external workbook link (xl/externalLinks/externalLink1.xml
# -> .../reg_results/multiple_loans.csv); that source file isn't reachable
# from this session, so the refreshed figures are written straight into the
# displayed table cells.
#
# NumberFormat is forced to Text ("@") before each write so that values such
# as "-38.9" or "0.030" are stored verbatim (keeping trailing zeros / minus
# signs) instead of being auto-coerced to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Row 7 / external row 8 & 14 (point estimates, top line of each panel) ---
Set-TextValue "B7" "-204.0***"
Set-TextValue "E7" "-0.066***"
Set-TextValue "G7" "-176.8***"
Set-TextValue "H7" "-0.080***"
Set-TextValue "J7" "-0.033"
Set-TextValue "L7" "-161.5***"

# --- Row 8 / external row 9 (standard errors) ---
Set-TextValue "B8" "(48.1)"

# --- Row 9 / external row 11 & 17 (second coefficient block) ---
Set-TextValue "B9" "-38.9"
Set-TextValue "C9" "-0.0086"
Set-TextValue "G9" "-33.3"
Set-TextValue "H9" "0.0037"
Set-TextValue "J9" "-0.0042"
Set-TextValue "L9" "-32.3"
Set-TextValue "M9" "-0.017"
Set-TextValue "O9" "-0.031*"

# --- Row 10 / external row 12 & 18 (standard errors) ---
Set-TextValue "G10" "(43.9)"
Set-TextValue "L10" "(40.6)"

# --- Row 13 / external row 21 (R-squared) ---
Set-TextValue "H13" "0.031"
Set-TextValue "M13" "0.030"

# --- Row 14 / external row 22 (control mean) ---
Set-TextValue "B14" "942.4"
Set-TextValue "E14" "0.44"
Set-TextValue "G14" "907.9"
Set-TextValue "L14" "907.9"
